# Generate Report for handoff
# This script swaps the display order / row content for the two tracked
# localization files (5d56e154-... and f646a3c9-...) across all three
# sheets, and marks the 5d56e154 file as "Ready for handoff" with a new
# handoff timestamp, reflecting a fresh handoff report being generated.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: "Overview"
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)

$ws1.Range("A2").Value = "f646a3c9-4bee-4f3d-ab89-8d7e7aeb7c52.md"
$ws1.Range("B2").Value = "Handed back: in sync with en-US"
$ws1.Range("C2").Value = "Handed back: in sync with en-US"

$ws1.Range("A3").Value = "5d56e154-5817-4996-8ea4-a2ce60afe0de.md"
$ws1.Range("B3").Value = "Ready for handoff"
$ws1.Range("C3").Value = "Ready for handoff"

$ws1.Range("A4").Value = ".localization-config"
$ws1.Range("B4").Value = "Not to be localized"
$ws1.Range("C4").Value = "Not to be localized"

# Rebuild the hyperlinks so the visible (display) text tracks the new
# cell contents while keeping the same external targets the original
# relationships pointed at.
$ws1.Hyperlinks.Delete()
$ws1.Hyperlinks.Add($ws1.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/5a5d6e297966b43d7e920a6522fb7087b7321a4f/e2e/5d56e154-5817-4996-8ea4-a2ce60afe0de.md", "", "", "f646a3c9-4bee-4f3d-ab89-8d7e7aeb7c52.md") | Out-Null
$ws1.Hyperlinks.Add($ws1.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/5a5d6e297966b43d7e920a6522fb7087b7321a4f/e2e/f646a3c9-4bee-4f3d-ab89-8d7e7aeb7c52.md", "", "", "5d56e154-5817-4996-8ea4-a2ce60afe0de.md") | Out-Null
$ws1.Hyperlinks.Add($ws1.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/5a5d6e297966b43d7e920a6522fb7087b7321a4f/.localization-config", "", "", ".localization-config") | Out-Null

# ---------------------------------------------------------------------
# Sheet 2: "zh-cn"
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)

$ws2.Range("A2").Value = "f646a3c9-4bee-4f3d-ab89-8d7e7aeb7c52.md"
$ws2.Range("B2").Value = "Handed back: in sync with en-US"
$ws2.Range("C2").Value = "f646a3c9-4bee-4f3d-ab89-8d7e7aeb7c52.011249e22c675102a82d84a292eef65c3d2d6459.zh-cn.xlf"
$ws2.Range("D2").Value = "2016-01-15 15:55:53"
$ws2.Range("E2").Value = "f646a3c9-4bee-4f3d-ab89-8d7e7aeb7c52.md"
$ws2.Range("F2").Value = "f646a3c9-4bee-4f3d-ab89-8d7e7aeb7c52.011249e22c675102a82d84a292eef65c3d2d6459.zh-cn.xlf"
$ws2.Range("G2").Value = "2016-01-15 15:56:58"
$ws2.Range("H2").Value = "Include"

$ws2.Range("A3").Value = "5d56e154-5817-4996-8ea4-a2ce60afe0de.md"
$ws2.Range("B3").Value = "Ready for handoff"
$ws2.Range("C3").Value = "5d56e154-5817-4996-8ea4-a2ce60afe0de.04e9c9b19b54a60110662a0b05de9adc835b9762.zh-cn.xlf"
$ws2.Range("D3").Value = "2016-01-15 15:57:57"
$ws2.Range("E3").Value = "5d56e154-5817-4996-8ea4-a2ce60afe0de.md"
$ws2.Range("F3").Value = "5d56e154-5817-4996-8ea4-a2ce60afe0de.04e9c9b19b54a60110662a0b05de9adc835b9762.zh-cn.xlf"
$ws2.Range("G3").Value = "2016-01-15 15:56:58"
$ws2.Range("H3").Value = "Include"

$ws2.Range("A4").Value = ".localization-config"
$ws2.Range("B4").Value = "Not to be localized"
$ws2.Range("D4").Value = "0001-01-01 00:00:00"
$ws2.Range("G4").Value = "0001-01-01 00:00:00"
$ws2.Range("H4").Value = "Ignored"

$ws2.Hyperlinks.Delete()
$ws2.Hyperlinks.Add($ws2.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/5a5d6e297966b43d7e920a6522fb7087b7321a4f/e2e/5d56e154-5817-4996-8ea4-a2ce60afe0de.md", "", "", "f646a3c9-4bee-4f3d-ab89-8d7e7aeb7c52.md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/71f115c370b597d65d7deacdb136ae5912183a24/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/yuwzho/5d56e154-5817-4996-8ea4-a2ce60afe0de.04e9c9b19b54a60110662a0b05de9adc835b9762.zh-cn.xlf", "", "", "f646a3c9-4bee-4f3d-ab89-8d7e7aeb7c52.011249e22c675102a82d84a292eef65c3d2d6459.zh-cn.xlf") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("E2"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/f0db7dc1bac67c40674a0d7910fe49d693c18a39/e2e/5d56e154-5817-4996-8ea4-a2ce60afe0de.md", "", "", "f646a3c9-4bee-4f3d-ab89-8d7e7aeb7c52.md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/f71e1552661a04a0c0c53583a93264ebd5bc49fa/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/yuwzho/5d56e154-5817-4996-8ea4-a2ce60afe0de.04e9c9b19b54a60110662a0b05de9adc835b9762.zh-cn.xlf", "", "", "f646a3c9-4bee-4f3d-ab89-8d7e7aeb7c52.011249e22c675102a82d84a292eef65c3d2d6459.zh-cn.xlf") | Out-Null

$ws2.Hyperlinks.Add($ws2.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/5a5d6e297966b43d7e920a6522fb7087b7321a4f/e2e/f646a3c9-4bee-4f3d-ab89-8d7e7aeb7c52.md", "", "", "5d56e154-5817-4996-8ea4-a2ce60afe0de.md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("C3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/71f115c370b597d65d7deacdb136ae5912183a24/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/yuwzho/f646a3c9-4bee-4f3d-ab89-8d7e7aeb7c52.011249e22c675102a82d84a292eef65c3d2d6459.zh-cn.xlf", "", "", "5d56e154-5817-4996-8ea4-a2ce60afe0de.04e9c9b19b54a60110662a0b05de9adc835b9762.zh-cn.xlf") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("E3"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/f0db7dc1bac67c40674a0d7910fe49d693c18a39/e2e/f646a3c9-4bee-4f3d-ab89-8d7e7aeb7c52.md", "", "", "5d56e154-5817-4996-8ea4-a2ce60afe0de.md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("F3"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/f71e1552661a04a0c0c53583a93264ebd5bc49fa/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/yuwzho/f646a3c9-4bee-4f3d-ab89-8d7e7aeb7c52.011249e22c675102a82d84a292eef65c3d2d6459.zh-cn.xlf", "", "", "5d56e154-5817-4996-8ea4-a2ce60afe0de.04e9c9b19b54a60110662a0b05de9adc835b9762.zh-cn.xlf") | Out-Null

$ws2.Hyperlinks.Add($ws2.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/5a5d6e297966b43d7e920a6522fb7087b7321a4f/.localization-config", "", "", ".localization-config") | Out-Null

# ---------------------------------------------------------------------
# Sheet 3: "de-de"
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item(3)

$ws3.Range("A2").Value = "f646a3c9-4bee-4f3d-ab89-8d7e7aeb7c52.md"
$ws3.Range("B2").Value = "Handed back: in sync with en-US"
$ws3.Range("C2").Value = "f646a3c9-4bee-4f3d-ab89-8d7e7aeb7c52.011249e22c675102a82d84a292eef65c3d2d6459.de-de.xlf"
$ws3.Range("D2").Value = "2016-01-15 15:56:05"
$ws3.Range("E2").Value = "f646a3c9-4bee-4f3d-ab89-8d7e7aeb7c52.md"
$ws3.Range("F2").Value = "f646a3c9-4bee-4f3d-ab89-8d7e7aeb7c52.011249e22c675102a82d84a292eef65c3d2d6459.de-de.xlf"
$ws3.Range("G2").Value = "2016-01-15 15:57:18"
$ws3.Range("H2").Value = "Include"

$ws3.Range("A3").Value = "5d56e154-5817-4996-8ea4-a2ce60afe0de.md"
$ws3.Range("B3").Value = "Ready for handoff"
$ws3.Range("C3").Value = "5d56e154-5817-4996-8ea4-a2ce60afe0de.04e9c9b19b54a60110662a0b05de9adc835b9762.de-de.xlf"
$ws3.Range("D3").Value = "2016-01-15 15:58:08"
$ws3.Range("E3").Value = "5d56e154-5817-4996-8ea4-a2ce60afe0de.md"
$ws3.Range("F3").Value = "5d56e154-5817-4996-8ea4-a2ce60afe0de.04e9c9b19b54a60110662a0b05de9adc835b9762.de-de.xlf"
$ws3.Range("G3").Value = "2016-01-15 15:57:18"
$ws3.Range("H3").Value = "Include"

$ws3.Range("A4").Value = ".localization-config"
$ws3.Range("B4").Value = "Not to be localized"
$ws3.Range("D4").Value = "0001-01-01 00:00:00"
$ws3.Range("G4").Value = "0001-01-01 00:00:00"
$ws3.Range("H4").Value = "Ignored"

$ws3.Hyperlinks.Delete()
$ws3.Hyperlinks.Add($ws3.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/5a5d6e297966b43d7e920a6522fb7087b7321a4f/e2e/5d56e154-5817-4996-8ea4-a2ce60afe0de.md", "", "", "f646a3c9-4bee-4f3d-ab89-8d7e7aeb7c52.md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/eae2abf79d3eadd93769713e19d132e4686d1902/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/yuwzho/5d56e154-5817-4996-8ea4-a2ce60afe0de.04e9c9b19b54a60110662a0b05de9adc835b9762.de-de.xlf", "", "", "f646a3c9-4bee-4f3d-ab89-8d7e7aeb7c52.011249e22c675102a82d84a292eef65c3d2d6459.de-de.xlf") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("E2"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/9397588e09ece4d4cb8e62416767e80027e4c508/e2e/5d56e154-5817-4996-8ea4-a2ce60afe0de.md", "", "", "f646a3c9-4bee-4f3d-ab89-8d7e7aeb7c52.md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/5dcb75321612716097d28ecc6b3529df3543531b/ol-handback/OpenLocalizationTestOrg/oltest.de-de/yuwzho/5d56e154-5817-4996-8ea4-a2ce60afe0de.04e9c9b19b54a60110662a0b05de9adc835b9762.de-de.xlf", "", "", "f646a3c9-4bee-4f3d-ab89-8d7e7aeb7c52.011249e22c675102a82d84a292eef65c3d2d6459.de-de.xlf") | Out-Null

$ws3.Hyperlinks.Add($ws3.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/5a5d6e297966b43d7e920a6522fb7087b7321a4f/e2e/f646a3c9-4bee-4f3d-ab89-8d7e7aeb7c52.md", "", "", "5d56e154-5817-4996-8ea4-a2ce60afe0de.md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("C3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/eae2abf79d3eadd93769713e19d132e4686d1902/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/yuwzho/f646a3c9-4bee-4f3d-ab89-8d7e7aeb7c52.011249e22c675102a82d84a292eef65c3d2d6459.de-de.xlf", "", "", "5d56e154-5817-4996-8ea4-a2ce60afe0de.04e9c9b19b54a60110662a0b05de9adc835b9762.de-de.xlf") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("E3"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/9397588e09ece4d4cb8e62416767e80027e4c508/e2e/f646a3c9-4bee-4f3d-ab89-8d7e7aeb7c52.md", "", "", "5d56e154-5817-4996-8ea4-a2ce60afe0de.md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("F3"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/5dcb75321612716097d28ecc6b3529df3543531b/ol-handback/OpenLocalizationTestOrg/oltest.de-de/yuwzho/f646a3c9-4bee-4f3d-ab89-8d7e7aeb7c52.011249e22c675102a82d84a292eef65c3d2d6459.de-de.xlf", "", "", "5d56e154-5817-4996-8ea4-a2ce60afe0de.04e9c9b19b54a60110662a0b05de9adc835b9762.de-de.xlf") | Out-Null

$ws3.Hyperlinks.Add($ws3.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/5a5d6e297966b43d7e920a6522fb7087b7321a4f/.localization-config", "", "", ".localization-config") | Out-Null
